# Update the LR-pairs sheet (Insl5-Rxfp3) with refreshed TPM-derived statistics.
# Existing rows 2-4 are updated in place and new rows 5-7 are appended so that
# every combination of sending cluster (FAPs, MuSCs) x target cluster (ECs, FAPs, MuSCs)
# is represented for the Insl5 -> Rxfp3 ligand-receptor pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Insl5"
$ws.Range("C2").Value = "Rxfp3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1188713333333333
$ws.Range("H2").Value = 0.356614
$ws.Range("I2").Value = 0.2549271348773238
$ws.Range("J2").Value = 0.339160140832479
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.014941
$ws.Range("N2").Value = 0.044823
$ws.Range("O2").Value = 0.02111135707456304
$ws.Range("P2").Value = 0.0256164071193072
$ws.Range("Q2").Value = 0.001776056591333333
$ws.Range("R2").Value = 0.015984509322
$ws.Range("S2").Value = 0.005381857772390477
$ws.Range("T2").Value = 0.008688064246206348

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Insl5"
$ws.Range("C3").Value = "Rxfp3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1188713333333333
$ws.Range("H3").Value = 0.356614
$ws.Range("I3").Value = 0.2549271348773238
$ws.Range("J3").Value = 0.339160140832479
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3193893333333333
$ws.Range("N3").Value = 0.9581679999999999
$ws.Range("O3").Value = 0.4512912296236288
$ws.Range("P3").Value = 0.54759435059439
$ws.Range("Q3").Value = 0.03796623590577777
$ws.Range("R3").Value = 0.3416961231519999
$ws.Range("S3").Value = 0.1150463801632161
$ws.Range("T3").Value = 0.1857221770666632

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Insl5"
$ws.Range("C4").Value = "Rxfp3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1188713333333333
$ws.Range("H4").Value = 0.356614
$ws.Range("I4").Value = 0.2549271348773238
$ws.Range("J4").Value = 0.339160140832479
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.373393
$ws.Range("N4").Value = 0.7467860000000001
$ws.Range("O4").Value = 0.5275974133018082
$ws.Range("P4").Value = 0.4267892422863028
$ws.Range("Q4").Value = 0.04438572376733334
$ws.Range("R4").Value = 0.266314342604
$ws.Range("S4").Value = 0.1344988969417172
$ws.Range("T4").Value = 0.1447498995196095

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Insl5"
$ws.Range("C5").Value = "Rxfp3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.347424
$ws.Range("H5").Value = 0.694848
$ws.Range("I5").Value = 0.7450728651226762
$ws.Range("J5").Value = 0.6608398591675211
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.014941
$ws.Range("N5").Value = 0.044823
$ws.Range("O5").Value = 0.02111135707456304
$ws.Range("P5").Value = 0.0256164071193072
$ws.Range("Q5").Value = 0.005190861984
$ws.Range("R5").Value = 0.031145171904
$ws.Range("S5").Value = 0.01572949930217257
$ws.Range("T5").Value = 0.01692834287310086

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Insl5"
$ws.Range("C6").Value = "Rxfp3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.347424
$ws.Range("H6").Value = 0.694848
$ws.Range("I6").Value = 0.7450728651226762
$ws.Range("J6").Value = 0.6608398591675211
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3193893333333333
$ws.Range("N6").Value = 0.9581679999999999
$ws.Range("O6").Value = 0.4512912296236288
$ws.Range("P6").Value = 0.54759435059439
$ws.Range("Q6").Value = 0.110963519744
$ws.Range("R6").Value = 0.665781118464
$ws.Range("S6").Value = 0.3362448494604127
$ws.Range("T6").Value = 0.3618721735277268

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Insl5"
$ws.Range("C7").Value = "Rxfp3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.347424
$ws.Range("H7").Value = 0.694848
$ws.Range("I7").Value = 0.7450728651226762
$ws.Range("J7").Value = 0.6608398591675211
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.373393
$ws.Range("N7").Value = 0.7467860000000001
$ws.Range("O7").Value = 0.5275974133018082
$ws.Range("P7").Value = 0.4267892422863028
$ws.Range("Q7").Value = 0.129725689632
$ws.Range("R7").Value = 0.5189027585280001
$ws.Range("S7").Value = 0.393098516360091
$ws.Range("T7").Value = 0.2820393427666934
